$wb = $excel.ActiveWorkbook

$wsConstrainedTemplate = $wb.Worksheets.Item("P1_Constrained")
$wsNotConstrainedTemplate = $wb.Worksheets.Item("P1_notConstrained")

# --- P2_Constrained (copy of P1_Constrained, shifted +10) ---
$wsConstrainedTemplate.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "P2_Constrained"

$ws3.Cells.Item(2,4).Value = 10
$ws3.Cells.Item(2,5).Value = 10
$ws3.Cells.Item(3,4).Value = 10
$ws3.Cells.Item(3,5).Value = 90
$ws3.Cells.Item(4,4).Value = 90
$ws3.Cells.Item(4,5).Value = 170
$ws3.Cells.Item(5,4).Value = 170
$ws3.Cells.Item(5,5).Value = 250
$ws3.Cells.Item(6,4).Value = 250
$ws3.Cells.Item(6,5).Value = 250

# --- P2_notConstrained (copy of P1_notConstrained, shifted +10) ---
$wsNotConstrainedTemplate.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "P2_notConstrained"

$ws4.Cells.Item(2,4).Value = 10
$ws4.Cells.Item(2,5).Value = 10
$ws4.Cells.Item(3,4).Value = 10
$ws4.Cells.Item(3,5).Value = 90
$ws4.Cells.Item(4,4).Value = 10
$ws4.Cells.Item(4,5).Value = 90
$ws4.Cells.Item(5,4).Value = 10
$ws4.Cells.Item(5,5).Value = 90
$ws4.Cells.Item(6,4).Value = 90
$ws4.Cells.Item(6,5).Value = 90

# --- P3_Constrained (copy of P1_Constrained, shifted +20) ---
$wsConstrainedTemplate.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5.Name = "P3_Constrained"

$ws5.Cells.Item(2,4).Value = 20
$ws5.Cells.Item(2,5).Value = 20
$ws5.Cells.Item(3,4).Value = 20
$ws5.Cells.Item(3,5).Value = 100
$ws5.Cells.Item(4,4).Value = 100
$ws5.Cells.Item(4,5).Value = 180
$ws5.Cells.Item(5,4).Value = 180
$ws5.Cells.Item(5,5).Value = 260
$ws5.Cells.Item(6,4).Value = 260
$ws5.Cells.Item(6,5).Value = 260

# --- P3_notConstrained (copy of P1_notConstrained, shifted +20) ---
$wsNotConstrainedTemplate.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6.Name = "P3_notConstrained"

$ws6.Cells.Item(2,4).Value = 20
$ws6.Cells.Item(2,5).Value = 20
$ws6.Cells.Item(3,4).Value = 20
$ws6.Cells.Item(3,5).Value = 100
$ws6.Cells.Item(4,4).Value = 20
$ws6.Cells.Item(4,5).Value = 100
$ws6.Cells.Item(5,4).Value = 20
$ws6.Cells.Item(5,5).Value = 100
$ws6.Cells.Item(6,4).Value = 100
$ws6.Cells.Item(6,5).Value = 100

$wsConstrainedTemplate.Select()
